# B6-PowerPoint.pptx edit
#
# 1) Every table in the deck currently uses the built-in "Table_0" table
#    style ({9F41B219-59DD-4D64-A3EF-F0F422A2D45A}); the commit re-styles
#    all of them to {B68B50B7-CD0A-4F65-B0BB-B95DD1D2121F}.
# 2) The presentation's theme (color scheme) is swapped from the
#    "Integral" / "Red Violet" palette to the stock Office "Office"
#    palette (font scheme / format scheme are identical between the two
#    themes already, only the 12 scheme colors - and the theme/clrScheme
#    display names, which aren't writable through this object model -
#    differ).

$p = $ppt.ActivePresentation

# --- 1. Re-style every table in the deck -----------------------------
$oldStyleId = "{9F41B219-59DD-4D64-A3EF-F0F422A2D45A}"
$newStyleId = "{B68B50B7-CD0A-4F65-B0BB-B95DD1D2121F}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2. Swap the theme color scheme onto the stock Office palette ----
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (RGB encoded the usual
# VBA way: R + G*256 + B*65536).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$design = $p.Designs.Item(1)
$scheme = $design.SlideMaster.Theme.ThemeColorScheme
for ($k = 1; $k -le $scheme.Count; $k++) {
    $scheme.Colors($k).RGB = $officeColors[$k - 1]
}
